# After splitting the test train method into multiple smaller modules,
# a new "Predicted_next_Day_Price" column (AB) was introduced and the
# downstream "Predicted_Signal" (AC) / "Actual_Return" (AD) values were
# recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column AB: Predicted_next_Day_Price values for rows 2-7 (all 0 in this run)
$ws.Range("AB2").Value = 0
$ws.Range("AB3").Value = 0
$ws.Range("AB4").Value = 0
$ws.Range("AB5").Value = 0
$ws.Range("AB6").Value = 0
$ws.Range("AB7").Value = 0

# Updated column AC: Predicted_Signal
$ws.Range("AC2").Value = 1
$ws.Range("AC3").Value = 1
$ws.Range("AC4").Value = 1
$ws.Range("AC5").Value = 0
$ws.Range("AC6").Value = 0
$ws.Range("AC7").Value = 0

# Updated column AD: Actual_Return
$ws.Range("AD3").Value = -0.01411042944785268
$ws.Range("AD4").Value = -0.009458618543870534
$ws.Range("AD5").Value = -0.08782510365623819
$ws.Range("AD6").Value = 0
$ws.Range("AD7").Value = 0
